$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2.177054233802296
$ws.Range("D2").Value = 2.94415561582861
$ws.Range("E2").Value = 16.66802848903046
$ws.Range("F2").Value = 24.47134540445481
$ws.Range("G2").Value = 3.567526063762902
$ws.Range("I2").Value = 19.59752123846195
$ws.Range("N2").Value = 18.99769123883711
$ws.Range("O2").Value = 20.89182271100775

$ws.Range("C3").Value = 2.172435147857365
$ws.Range("D3").Value = 2.95499900982323
$ws.Range("E3").Value = 15.70626053483035
$ws.Range("F3").Value = 23.77279097030611
$ws.Range("G3").Value = 3.571234858654964
$ws.Range("I3").Value = 19.17331663189001
$ws.Range("N3").Value = 18.40031508502701
$ws.Range("O3").Value = 20.41879351793042

$ws.Range("C4").Value = 2.1699388767347
$ws.Range("D4").Value = 2.962158412889302
$ws.Range("E4").Value = 15.09050342618864
$ws.Range("F4").Value = 23.34195226825164
$ws.Range("G4").Value = 3.573629020972537
$ws.Range("I4").Value = 18.91479532521317
$ws.Range("N4").Value = 18.02485520896362
$ws.Range("O4").Value = 20.12992364643306

$ws.Range("C5").Value = 2.169007673701893
$ws.Range("D5").Value = 2.965201582666245
$ws.Range("E5").Value = 14.8334996344148
$ws.Range("F5").Value = 23.16620146761137
$ws.Range("G5").Value = 3.574634190340252
$ws.Range("I5").Value = 18.81009776041359
$ws.Range("N5").Value = 17.86990355188765
$ws.Range("O5").Value = 20.01278767864523

$ws.Range("C6").Value = 2.168858261243953
$ws.Range("D6").Value = 2.965714475022706
$ws.Range("E6").Value = 14.79046609958031
$ws.Range("F6").Value = 23.13701640717197
$ws.Range("G6").Value = 3.574802884975969
$ws.Range("I6").Value = 18.79275700330097
$ws.Range("N6").Value = 17.84406337566597
$ws.Range("O6").Value = 19.99337790370629

$ws.Range("C7").Value = 2.16992596903508
$ws.Range("D7").Value = 2.96219894600785
$ws.Range("E7").Value = 15.08706160069387
$ws.Range("F7").Value = 23.33958235184754
$ws.Range("G7").Value = 3.573642457318429
$ws.Range("I7").Value = 18.91338048007347
$ws.Range("N7").Value = 18.02277304767603
$ws.Range("O7").Value = 20.12834131072626

$ws.Range("C8").Value = 2.175391400923772
$ws.Range("D8").Value = 2.947790048412414
$ws.Range("E8").Value = 16.34179060929392
$ws.Range("F8").Value = 24.23106040309856
$ws.Range("G8").Value = 3.568780657121915
$ws.Range("I8").Value = 19.45094988731143
$ws.Range("N8").Value = 18.79364780656867
$ws.Range("O8").Value = 20.72850622303068

$ws.Range("C9").Value = 2.188778008507189
$ws.Range("D9").Value = 2.923534966248327
$ws.Range("E9").Value = 18.73688618351982
$ws.Range("F9").Value = 25.95145186653697
$ws.Range("G9").Value = 3.56016906549611
$ws.Range("I9").Value = 20.51348828274695
$ws.Range("N9").Value = 20.2273683202997
$ws.Range("O9").Value = 21.90996168125022

$ws.Range("C10").Value = 2.200200076093394
$ws.Range("D10").Value = 2.908184729423248
$ws.Range("E10").Value = 20.4190547827733
$ws.Range("F10").Value = 27.18294755802058
$ws.Range("G10").Value = 3.554396683160884
$ws.Range("I10").Value = 21.29013090498111
$ws.Range("N10").Value = 21.22223697909767
$ws.Range("O10").Value = 22.77054280716412

$ws.Range("C11").Value = 2.205731232051607
$ws.Range("D11").Value = 2.901745097343441
$ws.Range("E11").Value = 21.14206437904654
$ws.Range("F11").Value = 27.73323559559251
$ws.Range("G11").Value = 3.551889431349226
$ws.Range("I11").Value = 21.64077404586619
$ws.Range("N11").Value = 21.66018057919901
$ws.Range("O11").Value = 23.15842820242483

$ws.Range("C12").Value = 2.207873030907942
$ws.Range("D12").Value = 2.89938535145791
$ws.Range("E12").Value = 21.40981913813115
$ws.Range("F12").Value = 27.93997367935093
$ws.Range("G12").Value = 3.5509569310822
$ws.Range("I12").Value = 21.77303194831448
$ws.Range("N12").Value = 21.82377585682186
$ws.Range("O12").Value = 23.30463903254842

$ws.Range("C13").Value = 2.207409669198068
$ws.Range("D13").Value = 2.899890048953178
$ws.Range("E13").Value = 21.35242109623703
$ws.Range("F13").Value = 27.89552529305429
$ws.Range("G13").Value = 3.551157010091716
$ws.Range("I13").Value = 21.74457326096319
$ws.Range("N13").Value = 21.78864458690801
$ws.Range("O13").Value = 23.27318221303199

$ws.Range("C14").Value = 2.205906497270679
$ws.Range("D14").Value = 2.901549376282606
$ws.Range("E14").Value = 21.16421342581512
$ws.Range("F14").Value = 27.75027818211127
$ws.Range("G14").Value = 3.551812375130608
$ws.Range("I14").Value = 21.65166624362908
$ws.Range("N14").Value = 21.67368539489659
$ws.Range("O14").Value = 23.17047139746499

$ws.Range("C15").Value = 2.204991888575577
$ws.Range("D15").Value = 2.902576045441832
$ws.Range("E15").Value = 21.04814619265927
$ws.Range("F15").Value = 27.66108993610336
$ws.Range("G15").Value = 3.552216007939526
$ws.Range("I15").Value = 21.59468585420258
$ws.Range("N15").Value = 21.60297336126124
$ws.Range("O15").Value = 23.1074659743189

$ws.Range("C16").Value = 2.199845268635152
$ws.Range("D16").Value = 2.908616564002029
$ws.Range("E16").Value = 20.37095809521568
$ws.Range("F16").Value = 27.14676622563257
$ws.Range("G16").Value = 3.55456291510376
$ws.Range("I16").Value = 21.26714985075624
$ws.Range("N16").Value = 21.19330956972086
$ws.Range("O16").Value = 22.74510765813716

$ws.Range("C17").Value = 2.196773189785098
$ws.Range("D17").Value = 2.912461845258352
$ws.Range("E17").Value = 19.94473900356374
$ws.Range("F17").Value = 26.82854251495638
$ws.Range("G17").Value = 3.556032967714352
$ws.Range("I17").Value = 21.06543382608167
$ws.Range("N17").Value = 20.93814219015166
$ws.Range("O17").Value = 22.52177723928153

$ws.Range("C18").Value = 2.195037820809328
$ws.Range("D18").Value = 2.914724641524562
$ws.Range("E18").Value = 19.69561312536645
$ws.Range("F18").Value = 26.64458609276432
$ws.Range("G18").Value = 3.556889676510139
$ws.Range("I18").Value = 20.9491700753218
$ws.Range("N18").Value = 20.79000725568362
$ws.Range("O18").Value = 22.39299357549152

$ws.Range("C19").Value = 2.194455711243231
$ws.Range("D19").Value = 2.915499540320131
$ws.Range("E19").Value = 19.61057936871486
$ws.Range("F19").Value = 26.58215009205571
$ws.Range("G19").Value = 3.557181666023975
$ws.Range("I19").Value = 20.90976804047995
$ws.Range("N19").Value = 20.73962067985785
$ws.Range("O19").Value = 22.34933788662787

$ws.Range("C20").Value = 2.197096952948375
$ws.Range("D20").Value = 2.912047215823479
$ws.Range("E20").Value = 19.99052197472236
$ws.Range("F20").Value = 26.86251508760304
$ws.Range("G20").Value = 3.555875322552067
$ws.Range("I20").Value = 21.08693294040393
$ws.Range("N20").Value = 20.96544799483449
$ws.Range("O20").Value = 22.54558644856834

$ws.Range("C21").Value = 2.206346740067187
$ws.Range("D21").Value = 2.901059847334219
$ws.Range("E21").Value = 21.21965802419746
$ws.Range("F21").Value = 27.79298702712013
$ws.Range("G21").Value = 3.551619419681133
$ws.Range("I21").Value = 21.67897056426095
$ws.Range("N21").Value = 21.70751365554064
$ws.Range("O21").Value = 23.2006595070544

$ws.Range("C22").Value = 2.212667072012168
$ws.Range("D22").Value = 2.894338630682314
$ws.Range("E22").Value = 21.9878201978674
$ws.Range("F22").Value = 28.39143486394593
$ws.Range("G22").Value = 3.548936635058988
$ws.Range("I22").Value = 22.06279337296474
$ws.Range("N22").Value = 22.17935961385674
$ws.Range("O22").Value = 23.62480026203102

$ws.Range("C23").Value = 2.209268943688253
$ws.Range("D23").Value = 2.897883577502701
$ws.Range("E23").Value = 21.58104139471089
$ws.Range("F23").Value = 28.07298321618894
$ws.Range("G23").Value = 3.550359496084504
$ws.Range("I23").Value = 21.85826868709075
$ws.Range("N23").Value = 21.92877110911181
$ws.Range("O23").Value = 23.39884195881073

$ws.Range("C24").Value = 2.196950483585789
$ws.Range("D24").Value = 2.912234507584127
$ws.Range("E24").Value = 19.96983622591278
$ws.Range("F24").Value = 26.847159202354
$ws.Range("G24").Value = 3.555946557926736
$ws.Range("I24").Value = 21.07721409737911
$ws.Range("N24").Value = 20.95310750188673
$ws.Range("O24").Value = 22.53482350175631

$ws.Range("C25").Value = 2.184874773439685
$ws.Range("D25").Value = 2.929665206824904
$ws.Range("E25").Value = 18.07972582065038
$ws.Range("F25").Value = 25.49069301138027
$ws.Range("G25").Value = 3.562400784297295
$ws.Range("I25").Value = 20.22611050585468
$ws.Range("N25").Value = 19.84905939529497
$ws.Range("O25").Value = 21.59094305220897
